$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")
$v = [string]$ws.Range("A8").Value
Write-Host $v
